$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ C = 2.319409367208825;   E = 5.080273296954374 }
    3  = @{ C = -3.942037578692481;  E = -1.648748515828502 }
    4  = @{ C = -2.839753013810475;  E = -1.632723506456935 }
    5  = @{ C = 4.960109259035406;   E = 4.506881698240095 }
    6  = @{ C = 2.134646894829806;   E = 2.531943146540772 }
    7  = @{ C = -2.700325749999488;  E = -0.3858735870725938 }
    8  = @{ C = 5.469647210234996;   E = 3.061326532789543 }
    9  = @{ C = 0.950153436409007;   E = 2.074800935750787 }
    10 = @{ C = 3.458696398997052;   E = 2.610227683091337 }
    11 = @{ C = 2.772413308959698;   E = 2.755099409670958 }
    12 = @{ C = 3.145819842658448;   E = 3.875106770584158 }
    13 = @{ C = 4.520465362328063;   E = 4.124307769579505 }
    14 = @{ C = 4.479055418855871;   E = 4.888255652935936 }
    15 = @{ C = 2.922623512367206;   E = 2.761298099516418 }
    16 = @{ C = -0.2939924376199055; E = 1.109922826784815 }
    17 = @{ C = -2.429883624035745;  E = -0.8739913853863412 }
    18 = @{ C = -1.060428249734879;  E = -1.285332554730334 }
    19 = @{ C = 0.1432745330888219;  E = -0.484131235569496 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
